$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: Humberto's record moves to period 2403 with updated values
$ws.Range("E16").Value = "2403"
$ws.Range("F16").Value = 52000
$ws.Range("G16").Value = 1300000

# Row 17: new record for Edelcy Cardossis Peña, period 2403
$ws.Range("C17").Value = "73092767"
$ws.Range("D17").Value = "EDELCY CARDOSSIS PEÑA"
$ws.Range("E17").Value = "2403"
$ws.Range("G17").Value = 1300000

# Row 18: Humberto's record for period 2404, restoring original F value
$ws.Range("C18").Value = "15122192"
$ws.Range("D18").Value = "HUMBERTO ANTONIO CALDERIN GONZALEZ"
$ws.Range("E18").Value = "2404"
$ws.Range("F18").Value = 1733
